$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28 (Tutorial_Eating / Bullet): mark as Fixed, clear the error message cell.
$ws.Range("B28").Value = "Fixed"
$ws.Range("D28").ClearContents()

# Row 29 (Tutorial_Eating / Vortex): mark as Fixed, clear the error message cell.
$ws.Range("B29").Value = "Fixed"
$ws.Range("D29").ClearContents()

# Update the view: scroll position and active selection moved down two rows.
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D28").Select()
